$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 0
    7 = 1
    8 = 2
    9 = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 3
    23 = 0
    24 = 0
    25 = 1
    26 = 3
    27 = 1
    28 = 0
    29 = 2
    30 = 2
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 2
    37 = 2
    38 = 2
    39 = 1
    40 = 2
    41 = 0
    42 = 2
    43 = 0
    44 = 3
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 1
    52 = 2
    53 = 0
    54 = 1
    56 = 1
    57 = 2
    58 = 4
    59 = 0
    60 = 1
    61 = 0
    62 = 2
    63 = 2
    64 = 1
    65 = 0
    66 = 1
    67 = 0
    69 = 1
    70 = 2
    71 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
